$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Russell Westbrook", "PG,SG", "Denver Nuggets")
    ,@("Dejounte Murray", "PG,SG", "New Orleans Pelicans")
    ,@("Keon Johnson", "PG,SG", "Brooklyn Nets")
    ,@("Chris Paul", "PG", "San Antonio Spurs")
    ,@("Jaylen Brown", "SG,SF", "Boston Celtics")
    ,@("Paolo Banchero", "SF,PF", "Orlando Magic")
    ,@("Pascal Siakam", "SF,PF,C", "Indiana Pacers")
    ,@("Buddy Hield", "SG,SF", "Golden State Warriors")
    ,@("Nikola Jokic", "C", "Denver Nuggets")
    ,@("Rudy Gobert", "C", "Minnesota Timberwolves")
    ,@("Jakob Poeltl", "C", "Toronto Raptors")
    ,@("Jalen Green", "PG,SG", "Houston Rockets")
    ,@("Bennedict Mathurin", "SG,SF", "Indiana Pacers")
    ,@("Aaron Wiggins", "SG,SF", "Oklahoma City Thunder")
    ,@("Deni Avdija", "SF,PF", "Portland Trail Blazers")
    ,@("Chet Holmgren", "PF,C", "Oklahoma City Thunder")
    ,@("Jalen Suggs", "PG,SG", "Orlando Magic")
    ,@("Khris Middleton", "SF", "Milwaukee Bucks")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
